# Paths in the Excel config are now relative to the Excel file itself
# (ampows adds its own directory to the path), so strip the leading
# "../examples/" segment from the generated sim/wind directories.

$wb = $excel.ActiveWorkbook

$configSheet = $wb.Worksheets.Item("config")

# Update the generated-output directory paths to be relative to the
# workbook location instead of climbing out to ../examples/.
$configSheet.Range("B10").Value = "./generated/sim"
$configSheet.Range("B11").Value = "./generated/wind"

# Move the active selection/cursor to B10 on the config sheet, and make
# the config sheet the active tab (it was DLC_List before).
$configSheet.Activate()
$configSheet.Range("B10").Select()
